$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.222.53"
$ws.Range("E2").Value = "  -3.71%  "

# Row 3
$ws.Range("D3").Value = "3.139.95"
$ws.Range("E3").Value = "  -3.30%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "601.82"
$ws.Range("E5").Value = "  -0.53%  "

# Row 6
$ws.Range("D6").Value = "145.95"
$ws.Range("E6").Value = "  -7.31%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").Value = "3.137.85"
$ws.Range("E8").Value = "  -3.33%  "

# Row 9
$ws.Range("D9").Value = "0.524"
$ws.Range("E9").Value = "  -4.52%  "

# Row 10
$ws.Range("E10").Value = "  -7.87%  "

# Row 11
$ws.Range("D11").Value = "5.48"
$ws.Range("E11").Value = "  -5.59%  "

# Row 12
$ws.Range("D12").Value = "0.472"
$ws.Range("E12").Value = "  -6.02%  "

# Row 13
$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  -8.69%  "

# Row 14
$ws.Range("D14").Value = "35.82"
$ws.Range("E14").Value = "  -8.41%  "

# Row 15
$ws.Range("D15").Value = "3.651.98"
$ws.Range("E15").Value = "  -3.19%  "

# Row 16
$ws.Range("D16").Value = "64.229.75"
$ws.Range("E16").Value = "  -3.68%  "

# Row 17
$ws.Range("E17").Value = "  +0.46%  "

# Row 18
$ws.Range("D18").Value = "3.137.95"
$ws.Range("E18").Value = "  -2.54%  "

# Row 19
$ws.Range("D19").Value = "6.90"
$ws.Range("E19").Value = "  -5.69%  "

# Row 20
$ws.Range("D20").Value = "476.64"
$ws.Range("E20").Value = "  -6.34%  "

# Row 21
$ws.Range("D21").Value = "14.57"
$ws.Range("E21").Value = "  -5.11%  "

# Row 22
$ws.Range("D22").Value = "0.704"
$ws.Range("E22").Value = "  -5.51%  "

# Row 23
$ws.Range("D23").Value = "7.65"
$ws.Range("E23").Value = "  -5.20%  "

# Row 24
$ws.Range("D24").Value = "13.63"
$ws.Range("E24").Value = "  -7.08%  "

# Row 25
$ws.Range("D25").Value = "83.41"
$ws.Range("E25").Value = "  -3.25%  "

# Row 26
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.09%  "

# Row 27
$ws.Range("E27").Value = "  -5.27%  "

# Row 28
$ws.Range("D28").Value = "8.35"
$ws.Range("E28").Value = "  -8.23%  "

# Row 29
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -7.70%  "

# Row 30
$ws.Range("D30").Value = "6.68"
$ws.Range("E30").Value = "  -3.46%  "

# Row 31
$ws.Range("E31").Value = "  -36.01%  "

# Row 32
$ws.Range("E32").Value = "  +0.13%  "

# Row 33
$ws.Range("D33").Value = "2.72"
$ws.Range("E33").Value = "  -6.26%  "

# Row 34
$ws.Range("D34").Value = "26.02"
$ws.Range("E34").Value = "  -8.01%  "

# Row 35
$ws.Range("E35").Value = "  -5.14%  "

# Row 36
$ws.Range("D36").Value = "54.09"
$ws.Range("E36").Value = "  -2.28%  "

# Row 37
$ws.Range("D37").Value = "5.95"
$ws.Range("E37").Value = "  -6.59%  "

# Row 38
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0707"
$ws.Range("E38").Value = "  -12.37%  "

# Row 39
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "443.41"
$ws.Range("E39").Value = "  -10.25%  "

# Row 40
$ws.Range("D40").Value = "2.87"
$ws.Range("E40").Value = "  -12.80%  "

# Row 41
$ws.Range("D41").Value = "0.0393"
$ws.Range("E41").Value = "  -7.39%  "

# Row 42
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "8.39"
$ws.Range("E42").Value = "  -4.18%  "

# Row 43
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.118"
$ws.Range("E43").Value = "  -8.14%  "

# Row 44
$ws.Range("D44").Value = "2.821.62"
$ws.Range("E44").Value = "  -4.29%  "

# Row 45
$ws.Range("D45").Value = "0.265"
$ws.Range("E45").Value = "  -10.02%  "

# Row 46
$ws.Range("D46").Value = "2.25"
$ws.Range("E46").Value = "  -9.08%  "

# Row 47
$ws.Range("E47").Value = "  -0.06%  "

# Row 48
$ws.Range("D48").Value = "26.25"
$ws.Range("E48").Value = "  -7.27%  "

# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "0.113"
$ws.Range("E49").Value = "  -4.88%  "

# Row 50
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "2.28"
$ws.Range("E50").Value = "  -5.57%  "

# Row 51
$ws.Range("D51").Value = "117.28"
$ws.Range("E51").Value = "  -3.30%  "
